# "titles, defns, d180, fixed catvars, zip inputs"
#
# - update the iode_quality_flag definition text on ColumnHeaders
# - add a new "CategoricalVariables" sheet (in front of ColumnHeaders)
#   describing the codes used by the toi_source / iode_quality_flag
#   categorical variables
# - tidy up the sheet selections left over from editing

$wb = $excel.ActiveWorkbook

# --- 1. fix up the existing ColumnHeaders sheet -----------------------
# Update this BEFORE adding the new sheet so the shared-string table
# keeps "flag" -> "IODE Quality Flag primary level" at its existing slot
# instead of appending a new (duplicate-ish) string at the end.
$colHeaders = $wb.Worksheets.Item("ColumnHeaders")
$colHeaders.Range("B13").Value = "IODE Quality Flag primary level"

# --- 2. add the new CategoricalVariables sheet -------------------------
$catVars = $wb.Worksheets.Add()
$catVars.Name = "CategoricalVariables"

$rows = @(
  @("attributeName", "code", "definition"),
  @("toi_source", "toi_niskin", "sample bottle was filled from a Niskin bottle on CTD rosette"),
  @("toi_source", "toi_underway", "sample bottle was filled from the ship's underway system"),
  @("iode_quality_flag", 1, "good"),
  @("iode_quality_flag", 2, "quality not evaluated, not available or unknown"),
  @("iode_quality_flag", 3, "questionable/suspect"),
  @("iode_quality_flag", 4, "bad"),
  @("iode_quality_flag", 9, "missing data")
)

for ($i = 0; $i -lt $rows.Count; $i++) {
  $r = $i + 1
  $catVars.Cells.Item($r, 1).Value = $rows[$i][0]
  $catVars.Cells.Item($r, 2).Value = $rows[$i][1]
  $catVars.Cells.Item($r, 3).Value = $rows[$i][2]
}

# --- 3. leave the selections where the author left them ----------------
# Re-fetch ColumnHeaders: the reference grabbed before Worksheets.Add()
# above no longer tracks the live sheet once a sibling sheet is inserted.
$colHeaders = $wb.Worksheets.Item("ColumnHeaders")
$colHeaders.Select()
$colHeaders.Range("B15").Select()

$catVars.Select()
$catVars.Range("B32").Select()
